# This script applies the edits described in the commit:
# "updated obj repo for study files tab; updated icdc profile"
#
# It updates the query text in row 2 (B2/C2) for the ICDC profile query
# (adds a Cohort column, and switches the StatQuery to a program-scoped
# variant), and it adds a new row 5 for a 'StudyFilesTab' entry that reuses
# the FilesTab query text (column B) and the StatQuery text (column C).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Text blocks (kept verbatim, including embedded blank lines / tabs) ----

# New StatQuery text for row 2, column C (program-scoped counts query)
$statQueryProgram = @'
MATCH (s:study)
  MATCH (demo:demographic) 
  MATCH (diag:diagnosis)
 MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis)
	WHERE s.clinical_study_designation IN ['MGT01'] and demo.breed in ['Australian Cattle Dog','Mixed Breed']and diag.disease_term in ['Mammary Cancer'] and diag.primary_disease_site in ['Mammary Gland']
OPTIONAL MATCH (s)<-[:member_of]-(c:case)
OPTIONAL MATCH (c)<-[:of_case]-(samp:sample)<-[:of_sample]-(f:file)
RETURN 
	count(DISTINCT(f)) as number_of_files , 
	count(DISTINCT(samp)) as number_of_sample , 
	count(DISTINCT(c.case_id)) as number_of_cases , 
	count(DISTINCT(s.clinical_study_designation)) as number_of_study
'@

# New detailed case query for row 2, column B (adds Cohort column)
$caseQueryWithCohort = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)

MATCH (c)<--(diag:diagnosis)
WHERE s.clinical_study_designation IN ['MGT01'] and demo.breed in ['Australian Cattle Dog','Mixed Breed']and diag.disease_term in ['Mammary Cancer'] and diag.primary_disease_site in ['Mammary Gland']
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment` 
        coalesce(co.cohort_description, '') AS `Cohort`
'@

# Existing Files query text, reused for the new StudyFilesTab row (column B)
$filesQuery = @'
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
 MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
WHERE s.clinical_study_designation IN ['MGT01'] and demo.breed in ['Australian Cattle Dog','Mixed Breed']and diag.disease_term in ['Mammary Cancer'] and diag.primary_disease_site in ['Mammary Gland']
WITH DISTINCT f, parent, c, demo, diag, s
RETURN coalesce(f.file_name, '') AS `File Name`, 
        coalesce(f.file_type, '') AS `File Type`, 
        coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(f.file_format, '') AS `File Format`,
        coalesce(f.file_size, '') AS `Size`,
        coalesce(c.case_id, '') AS `Case ID`, 
        coalesce(demo.breed,'') AS Breed , 
        coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(s.clinical_study_designation,'') AS `Study Code`
'@

# Existing StatQuery text (unchanged), reused for column C in rows 3, 4 and 5
$statQuery = @'
MATCH (s:study)
  MATCH (demo:demographic) 
  MATCH (diag:diagnosis)
 MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis)
	WHERE s.clinical_study_designation IN ['MGT01'] and demo.breed in ['Australian Cattle Dog','Mixed Breed']and diag.disease_term in ['Mammary Cancer'] and diag.primary_disease_site in ['Mammary Gland']
    
OPTIONAL MATCH (s)<-[:member_of]-(c:case)
OPTIONAL MATCH (c)<-[:of_case]-(samp:sample)<-[:of_sample]-(f:file)
RETURN 
	count(DISTINCT(f)) as number_of_files , 
	count(DISTINCT(samp)) as number_of_sample , 
	count(DISTINCT(c.case_id)) as number_of_cases , 
	count(DISTINCT(s.clinical_study_designation)) as number_of_study
'@

# File name columns (D = Neo4j data file, E = Web data file)
$neo4jFileName = @'
TC01_Canine_StudyMGT-Breed_Diagnosis_PrimDiseaseSite_Neo4jData.xlsx
'@
$webFileName = @'
TC01_Canine_StudyMGT-Breed_Diagnosis_PrimDiseaseSite_WebData.xlsx
'@

# ---- Update row 2 (existing CasesTab row) ----
$ws.Range("B2").Value = $caseQueryWithCohort
$ws.Range("C2").Value = $statQueryProgram
$ws.Rows.Item(2).RowHeight = 304.5

# ---- Add new row 5 (StudyFilesTab) ----
$ws.Range("A5").Value = "StudyFilesTab"
$ws.Range("B5").Value = $filesQuery
$ws.Range("B5").WrapText = $true
$ws.Range("C5").Value = $statQuery
$ws.Range("C5").WrapText = $true
$ws.Range("D5").Value = $neo4jFileName
$ws.Range("E5").Value = $webFileName
$ws.Rows.Item(5).RowHeight = 261

# ---- Column widths (closest achievable values to the target widths) ----
$ws.Columns.Item(1).ColumnWidth = 11.333333333333334
$ws.Columns.Item(5).ColumnWidth = 62.666666666666664

# ---- Update the active selection to match the edited workbook ----
$ws.Range("C5").Select()
